$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9479566812515259
$ws.Range("B1").Value = 1.480408549308777
$ws.Range("C1").Value = 3.545056104660034
$ws.Range("D1").Value = 3.074406862258911
$ws.Range("E1").Value = 1.5489661693573
